# Generate Report for Handback
# Refresh the "Correspond Handoff Datetime" (H) and
# "Correspond Handback DateTime" (K) timestamps for the
# 6ec12533-a576-4f89-80bc-5822835220eb source file row (row 2) on both the
# "zh-cn" and "de-de" handback-status sheets.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Sheets.Item("zh-cn")
$zhcn.Range("H2").Value = "2016-08-13 14:59:11"
$zhcn.Range("K2").Value = "2016-08-13 14:59:39"

$dede = $wb.Sheets.Item("de-de")
$dede.Range("H2").Value = "2016-08-13 14:59:19"
$dede.Range("K2").Value = "2016-08-13 14:59:48"
